$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Publisher value (row 9, column B)
$ws.Range("B9").Value = "Independent Trusted Third Party of the University Medicine Greifswald"

# Contact value (row 10, column B)
$ws.Range("B10").Value = "Independent Trusted Third Party of the University Medicine Greifswald (https://www.ths-greifswald.de/)"

# Description value (row 12, column B) - was previously empty
$ws.Range("B12").Value = "Expandable specification of the application context and/or scope of a consent policy or module "
